$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, shifting rows 6:107 down to 7:108.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly observation.
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44921
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 100112022
$ws.Cells.Item(6, 7).Value = "Arveja Verde"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 55
$ws.Cells.Item(6, 11).Value = 25000
$ws.Cells.Item(6, 12).Value = 25000
$ws.Cells.Item(6, 13).Value = 25000
$ws.Cells.Item(6, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(6, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(6, 16).Value = 1000
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
